$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update changed cells in existing rows (2-98) ---
$ws.Range("D2").Value = 44468
$ws.Range("J2").Value = 500
$ws.Range("K2").Value = 23000
$ws.Range("L2").Value = 25000
$ws.Range("M2").Value = 24000
$ws.Range("P2").Value = 960
$ws.Range("D3").Value = 45134
$ws.Range("J3").Value = 600
$ws.Range("K3").Value = 23000
$ws.Range("L3").Value = 25000
$ws.Range("M3").Value = 24000
$ws.Range("P3").Value = 960
$ws.Range("D4").Value = 44672
$ws.Range("H4").Value = "Sin especificar"
$ws.Range("J4").Value = 160
$ws.Range("K4").Value = 23000
$ws.Range("L4").Value = 25000
$ws.Range("M4").Value = 24000
$ws.Range("P4").Value = 960
$ws.Range("D5").Value = 44356
$ws.Range("J5").Value = 300
$ws.Range("K5").Value = 26000
$ws.Range("L5").Value = 28000
$ws.Range("M5").Value = 27000
$ws.Range("P5").Value = 1080
$ws.Range("D6").Value = 44671
$ws.Range("H6").Value = "Sin especificar"
$ws.Range("K6").Value = 23000
$ws.Range("L6").Value = 25000
$ws.Range("M6").Value = 24000
$ws.Range("P6").Value = 960
$ws.Range("D7").Value = 44384
$ws.Range("J7").Value = 400
$ws.Range("K7").Value = 26000
$ws.Range("L7").Value = 28000
$ws.Range("M7").Value = 27000
$ws.Range("P7").Value = 1080
$ws.Range("D8").Value = 44791
$ws.Range("H8").Value = "Perfection"
$ws.Range("J8").Value = 500
$ws.Range("K8").Value = 27000
$ws.Range("L8").Value = 29000
$ws.Range("M8").Value = 28000
$ws.Range("P8").Value = 1120
$ws.Range("D9").Value = 44818
$ws.Range("J9").Value = 400
$ws.Range("K9").Value = 24000
$ws.Range("L9").Value = 27000
$ws.Range("M9").Value = 25500
$ws.Range("P9").Value = 1020
$ws.Range("D10").Value = 44714
$ws.Range("J10").Value = 240
$ws.Range("K10").Value = 27000
$ws.Range("L10").Value = 28000
$ws.Range("M10").Value = 27500
$ws.Range("P10").Value = 1100
$ws.Range("D11").Value = 44391
$ws.Range("J11").Value = 100
$ws.Range("K11").Value = 26000
$ws.Range("L11").Value = 28000
$ws.Range("M11").Value = 27000
$ws.Range("P11").Value = 1080
$ws.Range("D12").Value = 44455
$ws.Range("J12").Value = 800
$ws.Range("K12").Value = 28000
$ws.Range("L12").Value = 30000
$ws.Range("M12").Value = 29000
$ws.Range("P12").Value = 1160
$ws.Range("D13").Value = 44475
$ws.Range("J13").Value = 1000
$ws.Range("K13").Value = 22000
$ws.Range("L13").Value = 24000
$ws.Range("M13").Value = 23000
$ws.Range("P13").Value = 920
$ws.Range("D14").Value = 44727
$ws.Range("J14").Value = 160
$ws.Range("K14").Value = 28000
$ws.Range("L14").Value = 30000
$ws.Range("M14").Value = 29000
$ws.Range("P14").Value = 1160
$ws.Range("D15").Value = 45112
$ws.Range("J15").Value = 1000
$ws.Range("K15").Value = 25000
$ws.Range("L15").Value = 26000
$ws.Range("M15").Value = 25500
$ws.Range("P15").Value = 1020
$ws.Range("D16").Value = 44371
$ws.Range("J16").Value = 500
$ws.Range("D17").Value = 44769
$ws.Range("H17").Value = "Perfection"
$ws.Range("J17").Value = 500
$ws.Range("K17").Value = 30000
$ws.Range("L17").Value = 32000
$ws.Range("M17").Value = 31000
$ws.Range("P17").Value = 1240
$ws.Range("D18").Value = 44490
$ws.Range("J18").Value = 500
$ws.Range("K18").Value = 16000
$ws.Range("L18").Value = 18000
$ws.Range("M18").Value = 17000
$ws.Range("P18").Value = 680
$ws.Range("D19").Value = 44377
$ws.Range("J19").Value = 500
$ws.Range("K19").Value = 26000
$ws.Range("L19").Value = 28000
$ws.Range("M19").Value = 27000
$ws.Range("P19").Value = 1080
$ws.Range("D20").Value = 44798
$ws.Range("J20").Value = 400
$ws.Range("K20").Value = 30000
$ws.Range("L20").Value = 32000
$ws.Range("M20").Value = 31000
$ws.Range("P20").Value = 1240
$ws.Range("D21").Value = 44826
$ws.Range("J21").Value = 520
$ws.Range("K21").Value = 28000
$ws.Range("L21").Value = 30000
$ws.Range("M21").Value = 29000
$ws.Range("P21").Value = 1160
$ws.Range("D22").Value = 44868
$ws.Range("J22").Value = 300
$ws.Range("K22").Value = 13000
$ws.Range("L22").Value = 15000
$ws.Range("M22").Value = 14000
$ws.Range("P22").Value = 560
$ws.Range("D23").Value = 44763
$ws.Range("J23").Value = 400
$ws.Range("K23").Value = 29000
$ws.Range("L23").Value = 30000
$ws.Range("M23").Value = 29500
$ws.Range("P23").Value = 1180
$ws.Range("D24").Value = 44755
$ws.Range("J24").Value = 200
$ws.Range("K24").Value = 30000
$ws.Range("L24").Value = 32000
$ws.Range("M24").Value = 31000
$ws.Range("P24").Value = 1240
$ws.Range("D25").Value = 44749
$ws.Range("J25").Value = 470
$ws.Range("K25").Value = 28000
$ws.Range("L25").Value = 30000
$ws.Range("M25").Value = 29064
$ws.Range("P25").Value = 1163
$ws.Range("D26").Value = 44839
$ws.Range("J26").Value = 700
$ws.Range("K26").Value = 22000
$ws.Range("L26").Value = 24000
$ws.Range("M26").Value = 23000
$ws.Range("P26").Value = 920
$ws.Range("D27").Value = 44357
$ws.Range("J27").Value = 340
$ws.Range("K27").Value = 28000
$ws.Range("L27").Value = 30000
$ws.Range("M27").Value = 29000
$ws.Range("P27").Value = 1160
$ws.Range("D28").Value = 44349
$ws.Range("J28").Value = 600
$ws.Range("D29").Value = 44721
$ws.Range("J29").Value = 240
$ws.Range("K29").Value = 28000
$ws.Range("L29").Value = 30000
$ws.Range("M29").Value = 29000
$ws.Range("P29").Value = 1160
$ws.Range("D30").Value = 44804
$ws.Range("J30").Value = 400
$ws.Range("D31").Value = 44748
$ws.Range("J31").Value = 700
$ws.Range("K31").Value = 28000
$ws.Range("L31").Value = 30000
$ws.Range("M31").Value = 29000
$ws.Range("P31").Value = 1160
$ws.Range("D32").Value = 44413
$ws.Range("J32").Value = 700
$ws.Range("K32").Value = 26000
$ws.Range("L32").Value = 28000
$ws.Range("M32").Value = 27000
$ws.Range("P32").Value = 1080
$ws.Range("D33").Value = 44784
$ws.Range("J33").Value = 360
$ws.Range("K33").Value = 27000
$ws.Range("L33").Value = 29000
$ws.Range("M33").Value = 28000
$ws.Range("P33").Value = 1120
$ws.Range("D34").Value = 45092
$ws.Range("J34").Value = 300
$ws.Range("K34").Value = 27000
$ws.Range("L34").Value = 29000
$ws.Range("M34").Value = 28000
$ws.Range("P34").Value = 1120
$ws.Range("D35").Value = 44483
$ws.Range("K35").Value = 18000
$ws.Range("L35").Value = 20000
$ws.Range("M35").Value = 19000
$ws.Range("P35").Value = 760
$ws.Range("D36").Value = 44874
$ws.Range("H36").Value = "Perfection"
$ws.Range("J36").Value = 160
$ws.Range("K36").Value = 14000
$ws.Range("L36").Value = 16000
$ws.Range("M36").Value = 15000
$ws.Range("P36").Value = 600
$ws.Range("D37").Value = 45126
$ws.Range("J37").Value = 600
$ws.Range("K37").Value = 22000
$ws.Range("L37").Value = 24000
$ws.Range("M37").Value = 23000
$ws.Range("P37").Value = 920
$ws.Range("D38").Value = 44679
$ws.Range("J38").Value = 400
$ws.Range("K38").Value = 25000
$ws.Range("L38").Value = 27000
$ws.Range("M38").Value = 26000
$ws.Range("P38").Value = 1040
$ws.Range("D39").Value = 44742
$ws.Range("J39").Value = 200
$ws.Range("K39").Value = 28000
$ws.Range("L39").Value = 30000
$ws.Range("M39").Value = 29000
$ws.Range("P39").Value = 1160
$ws.Range("D40").Value = 44461
$ws.Range("K40").Value = 23000
$ws.Range("L40").Value = 25000
$ws.Range("M40").Value = 24000
$ws.Range("P40").Value = 960
$ws.Range("D41").Value = 44392
$ws.Range("J41").Value = 100
$ws.Range("K41").Value = 26000
$ws.Range("L41").Value = 28000
$ws.Range("M41").Value = 27000
$ws.Range("P41").Value = 1080
$ws.Range("D42").Value = 44847
$ws.Range("H42").Value = "Sin especificar"
$ws.Range("J42").Value = 300
$ws.Range("K42").Value = 23000
$ws.Range("L42").Value = 24000
$ws.Range("M42").Value = 23500
$ws.Range("P42").Value = 940
$ws.Range("D43").Value = 44385
$ws.Range("J43").Value = 500
$ws.Range("K43").Value = 26000
$ws.Range("L43").Value = 28000
$ws.Range("M43").Value = 27000
$ws.Range("P43").Value = 1080
$ws.Range("D44").Value = 44783
$ws.Range("J44").Value = 400
$ws.Range("K44").Value = 27000
$ws.Range("L44").Value = 29000
$ws.Range("M44").Value = 28000
$ws.Range("P44").Value = 1120
$ws.Range("D45").Value = 44476
$ws.Range("K45").Value = 23000
$ws.Range("L45").Value = 24000
$ws.Range("M45").Value = 23500
$ws.Range("P45").Value = 940
$ws.Range("D46").Value = 44433
$ws.Range("J46").Value = 400
$ws.Range("K46").Value = 28000
$ws.Range("L46").Value = 30000
$ws.Range("M46").Value = 29000
$ws.Range("P46").Value = 1160
$ws.Range("D47").Value = 44707
$ws.Range("J47").Value = 200
$ws.Range("K47").Value = 25000
$ws.Range("L47").Value = 27000
$ws.Range("M47").Value = 26000
$ws.Range("P47").Value = 1040
$ws.Range("D48").Value = 45113
$ws.Range("J48").Value = 800
$ws.Range("K48").Value = 23000
$ws.Range("L48").Value = 25000
$ws.Range("M48").Value = 24000
$ws.Range("P48").Value = 960
$ws.Range("D49").Value = 44419
$ws.Range("J49").Value = 600
$ws.Range("K49").Value = 27000
$ws.Range("L49").Value = 29000
$ws.Range("M49").Value = 28000
$ws.Range("P49").Value = 1120
$ws.Range("D50").Value = 44812
$ws.Range("K50").Value = 28000
$ws.Range("L50").Value = 30000
$ws.Range("M50").Value = 29000
$ws.Range("P50").Value = 1160
$ws.Range("D51").Value = 44427
$ws.Range("K51").Value = 28000
$ws.Range("L51").Value = 30000
$ws.Range("M51").Value = 29000
$ws.Range("P51").Value = 1160
$ws.Range("D52").Value = 45085
$ws.Range("K52").Value = 25000
$ws.Range("L52").Value = 27000
$ws.Range("M52").Value = 26000
$ws.Range("P52").Value = 1040
$ws.Range("D53").Value = 44678
$ws.Range("H53").Value = "Perfection"
$ws.Range("J53").Value = 600
$ws.Range("K53").Value = 25000
$ws.Range("L53").Value = 27000
$ws.Range("M53").Value = 26000
$ws.Range("P53").Value = 1040
$ws.Range("D54").Value = 44448
$ws.Range("K54").Value = 28000
$ws.Range("L54").Value = 30000
$ws.Range("M54").Value = 29000
$ws.Range("P54").Value = 1160
$ws.Range("D55").Value = 44350
$ws.Range("J55").Value = 700
$ws.Range("K55").Value = 28000
$ws.Range("M55").Value = 29000
$ws.Range("P55").Value = 1160
$ws.Range("D56").Value = 44434
$ws.Range("J56").Value = 500
$ws.Range("K56").Value = 28000
$ws.Range("L56").Value = 30000
$ws.Range("M56").Value = 29000
$ws.Range("P56").Value = 1160
$ws.Range("D58").Value = 45120
$ws.Range("J58").Value = 1100
$ws.Range("K58").Value = 25000
$ws.Range("L58").Value = 27000
$ws.Range("M58").Value = 26000
$ws.Range("P58").Value = 1040
$ws.Range("D59").Value = 44441
$ws.Range("J59").Value = 700
$ws.Range("K59").Value = 28000
$ws.Range("L59").Value = 30000
$ws.Range("M59").Value = 29000
$ws.Range("P59").Value = 1160
$ws.Range("D60").Value = 45127
$ws.Range("J60").Value = 700
$ws.Range("K60").Value = 20000
$ws.Range("L60").Value = 22000
$ws.Range("M60").Value = 21000
$ws.Range("P60").Value = 840
$ws.Range("D61").Value = 44412
$ws.Range("J61").Value = 600
$ws.Range("K61").Value = 25000
$ws.Range("L61").Value = 27000
$ws.Range("M61").Value = 26000
$ws.Range("P61").Value = 1040
$ws.Range("D62").Value = 44832
$ws.Range("J62").Value = 600
$ws.Range("K62").Value = 23000
$ws.Range("L62").Value = 25000
$ws.Range("M62").Value = 24000
$ws.Range("P62").Value = 960
$ws.Range("D63").Value = 44406
$ws.Range("J63").Value = 600
$ws.Range("K63").Value = 26000
$ws.Range("L63").Value = 28000
$ws.Range("M63").Value = 27000
$ws.Range("P63").Value = 1080
$ws.Range("D64").Value = 44776
$ws.Range("J64").Value = 400
$ws.Range("D65").Value = 44706
$ws.Range("J65").Value = 160
$ws.Range("K65").Value = 25000
$ws.Range("L65").Value = 26000
$ws.Range("M65").Value = 25500
$ws.Range("P65").Value = 1020
$ws.Range("D66").Value = 44426
$ws.Range("D67").Value = 44399
$ws.Range("J67").Value = 400
$ws.Range("D68").Value = 44685
$ws.Range("J68").Value = 160
$ws.Range("K68").Value = 25000
$ws.Range("L68").Value = 27000
$ws.Range("M68").Value = 26000
$ws.Range("P68").Value = 1040
$ws.Range("D69").Value = 44363
$ws.Range("J69").Value = 240
$ws.Range("D70").Value = 45084
$ws.Range("J70").Value = 500
$ws.Range("D71").Value = 44811
$ws.Range("D72").Value = 44462
$ws.Range("J72").Value = 400
$ws.Range("K72").Value = 22000
$ws.Range("L72").Value = 23000
$ws.Range("M72").Value = 22500
$ws.Range("P72").Value = 900
$ws.Range("D73").Value = 44699
$ws.Range("J73").Value = 200
$ws.Range("K73").Value = 29000
$ws.Range("M73").Value = 29500
$ws.Range("P73").Value = 1180
$ws.Range("D74").Value = 45091
$ws.Range("J74").Value = 360
$ws.Range("K74").Value = 26000
$ws.Range("L74").Value = 28000
$ws.Range("M74").Value = 27000
$ws.Range("P74").Value = 1080
$ws.Range("D75").Value = 44482
$ws.Range("K75").Value = 18000
$ws.Range("L75").Value = 20000
$ws.Range("M75").Value = 19000
$ws.Range("P75").Value = 760
$ws.Range("D76").Value = 44497
$ws.Range("J76").Value = 500
$ws.Range("K76").Value = 13000
$ws.Range("L76").Value = 15000
$ws.Range("M76").Value = 14000
$ws.Range("P76").Value = 560
$ws.Range("D77").Value = 44762
$ws.Range("J77").Value = 400
$ws.Range("K77").Value = 29000
$ws.Range("M77").Value = 29500
$ws.Range("P77").Value = 1180
$ws.Range("D78").Value = 44398
$ws.Range("J78").Value = 500
$ws.Range("D79").Value = 44435
$ws.Range("J79").Value = 900
$ws.Range("K79").Value = 28000
$ws.Range("L79").Value = 30000
$ws.Range("M79").Value = 29000
$ws.Range("P79").Value = 1160
$ws.Range("D80").Value = 44825
$ws.Range("J80").Value = 480
$ws.Range("K80").Value = 28000
$ws.Range("M80").Value = 29000
$ws.Range("P80").Value = 1160
$ws.Range("D81").Value = 44454
$ws.Range("J81").Value = 1000
$ws.Range("K81").Value = 28000
$ws.Range("L81").Value = 30000
$ws.Range("M81").Value = 29000
$ws.Range("P81").Value = 1160
$ws.Range("D82").Value = 44741
$ws.Range("J82").Value = 160
$ws.Range("D83").Value = 44790
$ws.Range("J83").Value = 560
$ws.Range("K83").Value = 27000
$ws.Range("L83").Value = 29000
$ws.Range("M83").Value = 28000
$ws.Range("P83").Value = 1120
$ws.Range("D84").Value = 44720
$ws.Range("J84").Value = 400
$ws.Range("K84").Value = 28000
$ws.Range("L84").Value = 30000
$ws.Range("M84").Value = 29000
$ws.Range("P84").Value = 1160
$ws.Range("D85").Value = 44370
$ws.Range("J85").Value = 400
$ws.Range("K85").Value = 27000
$ws.Range("L85").Value = 28000
$ws.Range("M85").Value = 27500
$ws.Range("P85").Value = 1100
$ws.Range("D86").Value = 44489
$ws.Range("J86").Value = 400
$ws.Range("K86").Value = 18000
$ws.Range("L86").Value = 20000
$ws.Range("M86").Value = 19000
$ws.Range("P86").Value = 760
$ws.Range("D87").Value = 44469
$ws.Range("J87").Value = 600
$ws.Range("K87").Value = 22000
$ws.Range("L87").Value = 24000
$ws.Range("M87").Value = 23000
$ws.Range("P87").Value = 920
$ws.Range("D88").Value = 45106
$ws.Range("J88").Value = 900
$ws.Range("K88").Value = 26000
$ws.Range("L88").Value = 28000
$ws.Range("M88").Value = 27000
$ws.Range("P88").Value = 1080
$ws.Range("D89").Value = 44860
$ws.Range("J89").Value = 200
$ws.Range("K89").Value = 15000
$ws.Range("L89").Value = 16000
$ws.Range("M89").Value = 15500
$ws.Range("P89").Value = 620
$ws.Range("D90").Value = 44447
$ws.Range("J90").Value = 600
$ws.Range("K90").Value = 28000
$ws.Range("L90").Value = 30000
$ws.Range("M90").Value = 29000
$ws.Range("P90").Value = 1160
$ws.Range("D91").Value = 44420
$ws.Range("J91").Value = 700
$ws.Range("K91").Value = 27000
$ws.Range("L91").Value = 29000
$ws.Range("M91").Value = 28000
$ws.Range("P91").Value = 1120
$ws.Range("D92").Value = 44819
$ws.Range("J92").Value = 500
$ws.Range("K92").Value = 25000
$ws.Range("L92").Value = 28000
$ws.Range("M92").Value = 26500
$ws.Range("P92").Value = 1060
$ws.Range("D93").Value = 45140
$ws.Range("J93").Value = 300
$ws.Range("K93").Value = 20000
$ws.Range("L93").Value = 22000
$ws.Range("M93").Value = 21000
$ws.Range("P93").Value = 840
$ws.Range("D94").Value = 45119
$ws.Range("J94").Value = 1000
$ws.Range("K94").Value = 26000
$ws.Range("L94").Value = 28000
$ws.Range("M94").Value = 27000
$ws.Range("P94").Value = 1080
$ws.Range("D95").Value = 44756
$ws.Range("J95").Value = 240
$ws.Range("K95").Value = 30000
$ws.Range("L95").Value = 32000
$ws.Range("M95").Value = 31000
$ws.Range("P95").Value = 1240
$ws.Range("D96").Value = 44343
$ws.Range("J96").Value = 200
$ws.Range("K96").Value = 26000
$ws.Range("L96").Value = 28000
$ws.Range("M96").Value = 27000
$ws.Range("P96").Value = 1080
$ws.Range("D97").Value = 44846
$ws.Range("H97").Value = "Sin especificar"
$ws.Range("J97").Value = 488
$ws.Range("K97").Value = 23000
$ws.Range("L97").Value = 24000
$ws.Range("M97").Value = 23426
$ws.Range("P97").Value = 937
$ws.Range("D98").Value = 45133
$ws.Range("J98").Value = 560
$ws.Range("K98").Value = 23000
$ws.Range("L98").Value = 25000
$ws.Range("M98").Value = 24000
$ws.Range("P98").Value = 960

# --- Append new row 99 (full record) ---
$ws.Range("A99").Value = 2
$ws.Range("B99").Value = "Comercializadora del Agro de Limarí"
$ws.Range("C99").Value = "Coquimbo"
$ws.Range("D99").Value = 44364
$ws.Range("E99").Value = 4
$ws.Range("F99").Value = 100112022
$ws.Range("G99").Value = "Arveja Verde"
$ws.Range("H99").Value = "Perfection"
$ws.Range("I99").Value = "Primera"
$ws.Range("J99").Value = 200
$ws.Range("K99").Value = 28000
$ws.Range("L99").Value = 30000
$ws.Range("M99").Value = 29000
$ws.Range("N99").Value = "`$/malla 25 kilos"
$ws.Range("O99").Value = "Provincia de Limarí"
$ws.Range("P99").Value = 1160
$ws.Range("Q99").Value = 25
$ws.Range("R99").Value = "Hortaliza"

# Match the date style/number format used by the other rows in column D
$ws.Range("D99").NumberFormat = $ws.Range("D98").NumberFormat

